# Append two new review rows (6 and 7) to the sheet, matching the
# formatting of the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: com.singleton.stretchy / taxi game review ---
$ws.Range("A6").Value = "com.singleton.stretchy"
$ws.Range("B6").Value = "taxi game"
$ws.Range("C6").Value = "nirh94846@gmail.com"
$ws.Range("D6").Value = "shamirnaftali@gmail.com"
$ws.Range("E6").Value = "27/5/2019 15:59"
$ws.Range("F6").Value = "be the best to make some plank between cities. I love this game so much. Great car game"

# --- Row 7: com.hamxa.shaynachim / bitcoin guide review ---
$ws.Range("A7").Value = "com.hamxa.shaynachim"
$ws.Range("B7").Value = "bitcoin guide"
$ws.Range("C7").Value = "nirh94846@gmail.com"
$ws.Range("D7").Value = "shamirnaftali@gmail.com"
$ws.Range("E7").Value = "27/5/2019 15:59"
$ws.Range("F7").Value = "very nice guide about bitcoin. I understood everything and now lets start working"

# Match the appid/review-text styling (col A and F) used by the rest of the
# table (Mangal font) ...
$ws.Range("A6").Font.Name = "Mangal"
$ws.Range("F6").Font.Name = "Mangal"
$ws.Range("A7").Font.Name = "Mangal"
$ws.Range("F7").Font.Name = "Mangal"

# ... and the email/recovery styling (col C and D: centered Calibri, black)
$ws.Range("C6:D6").Font.Name = "Calibri"
$ws.Range("C6:D6").Font.Size = 11
$ws.Range("C6:D6").Font.Color = 0
$ws.Range("C6:D6").HorizontalAlignment = -4108

$ws.Range("C7:D7").Font.Name = "Calibri"
$ws.Range("C7:D7").Font.Size = 11
$ws.Range("C7:D7").Font.Color = 0
$ws.Range("C7:D7").HorizontalAlignment = -4108

# Rows 4/5 (the other multi-line review rows) are slightly taller than the
# sheet default - match that for the two new rows too.
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Leave the selection on the last filled cell, like the source edit did.
$ws.Range("F7").Select() | Out-Null
